$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new test-case rows (111 and 112) describing the new
# "update scalar index filtered by its own value" cases ---
# Cell values are set in an order that mirrors how the new shared
# strings were introduced so the workbook content matches as closely
# as possible.

$ws.Range("A111").Value = "updel_110"
$ws.Range("A112").Value = "updel_111"

$ws.Range("F111").Value = "scalar055"
$ws.Range("G111").Value = "scalar055_value1"
$ws.Range("H111").Value = 'update $scalar055 set amount=888.88 where amount=279540.148'
$ws.Range("I111").Value = "1"
$ws.Range("J111").Value = 'select * from $scalar055 where amount=888.88'
$ws.Range("K111").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_110.csv"

$ws.Range("C111").Value = "更新标量索引字段值通过标量索引字段过滤"
$ws.Range("C112").Value = "更新标量索引字段值通过主键字段过滤"

$ws.Range("H112").Value = 'update $scalar055 set amount=-999.999 where id=4695'
$ws.Range("J112").Value = 'select * from $scalar055 where amount=-999.999'
$ws.Range("K112").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_111.csv"

$ws.Range("B111").Value = "y"
$ws.Range("D111").Value = "Index"
$ws.Range("E111").Value = "scalar_index"
$ws.Range("L111").Value = "csv_containsAll"

$ws.Range("B112").Value = "y"
$ws.Range("D112").Value = "Index"
$ws.Range("E112").Value = "scalar_index"
$ws.Range("F112").Value = "scalar055"
$ws.Range("G112").Value = "scalar055_value1"
$ws.Range("I112").Value = "1"
$ws.Range("L112").Value = "csv_containsAll"

# --- Widen column H to fit the new, longer SQL text ---
# (the platform rounds column width to whole-pixel increments, so this is
# the input that lands closest to the target display width of 60.625)
$ws.Columns.Item(8).ColumnWidth = 59.857142857142854

# --- Update the view state: scroll down and move the active selection ---
$ws.Range("G123").Select()
